$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.2248
$ws.Range("A3").Value = -21.8663
$ws.Range("D3").Value = -7.347699999999995
$ws.Range("E6").Value = 16.52210000000001
$ws.Range("D12").Value = -7.410000000000005
$ws.Range("A14").Value = -21.7693
$ws.Range("E19").Value = 16.46660000000001
$ws.Range("A21").Value = -19.94889999999998
$ws.Range("A23").Value = -20.47949999999998
$ws.Range("D24").Value = -7.519000000000003
$ws.Range("E24").Value = 17.26170000000001
$ws.Range("A25").Value = -22.00249999999999
$ws.Range("C25").Value = -13.1541
$ws.Range("D25").Value = -8.682099999999997
$ws.Range("A26").Value = -21.17679999999997
$ws.Range("C27").Value = -12.4281
$ws.Range("A29").Value = -20.95029999999998
$ws.Range("E30").Value = 15.43519999999999
$ws.Range("C31").Value = -12.8414
$ws.Range("E31").Value = 16.35950000000001
$ws.Range("E33").Value = 17.03130000000002
$ws.Range("C39").Value = -12.58070000000001
$ws.Range("E42").Value = 16.52650000000001
$ws.Range("C48").Value = -11.7187
$ws.Range("D50").Value = -8.065600000000003
$ws.Range("C51").Value = -11.6283
$ws.Range("C52").Value = -11.1229
$ws.Range("A53").Value = -22.28280000000001
$ws.Range("D53").Value = -6.3003
$ws.Range("C55").Value = -13.8203
$ws.Range("E55").Value = 16.46580000000001
$ws.Range("C56").Value = -11.7407
$ws.Range("A57").Value = -21.84230000000001
$ws.Range("C57").Value = -13.22709999999999
$ws.Range("D57").Value = -8.6881
$ws.Range("E58").Value = 16.58280000000001
$ws.Range("A59").Value = -22.2986
$ws.Range("D61").Value = -7.711299999999997
$ws.Range("D63").Value = -7.975300000000003
$ws.Range("E65").Value = 17.00280000000001
$ws.Range("A69").Value = -21.5737
$ws.Range("D70").Value = -8.124600000000004
$ws.Range("E70").Value = 16.7244
$ws.Range("C73").Value = -12.7918
$ws.Range("E75").Value = 16.56680000000001
$ws.Range("A79").Value = -20.3913
$ws.Range("A83").Value = -21.86719999999999
$ws.Range("E83").Value = 16.59950000000001
$ws.Range("D86").Value = -8.811000000000003
$ws.Range("E86").Value = 16.11570000000001
$ws.Range("C89").Value = -10.4396
$ws.Range("C90").Value = -12.5371
$ws.Range("A91").Value = -21.37360000000001
$ws.Range("C92").Value = -10.6222
$ws.Range("A93").Value = -21.01389999999998
$ws.Range("E96").Value = 15.68079999999999
$ws.Range("E97").Value = 17.00100000000002
$ws.Range("D98").Value = -9.056299999999995
$ws.Range("D100").Value = -8.600599999999996
$ws.Range("D102").Value = -7.710699999999998
